$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").NumberFormat = "m/d/yy h:mm"
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Range("A6").Value = $epoch.AddDays(42601.914861111109)

$ws.Range("B6").Value = "Gilead Sciences, Inc."
$ws.Range("C6").Value = "GILD"
$ws.Range("D6").Value = 80.91
$ws.Range("E6").Value = 80.7
$ws.Range("F6").Value = -0.14000000000000001
$ws.Range("G6").Value = 80.59
